# Update betting-odds / score values for the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 4.1
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 4.5
$ws.Range("AE4").Value = 17
$ws.Range("AL4").Value = 17
$ws.Range("AX4").Value = 11

# Row 5
$ws.Range("G5").Value = 1.55
$ws.Range("I5").Value = 6.5
$ws.Range("J5").Value = 2.2
$ws.Range("L5").Value = 6.5
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("X5").Value = 6.5
$ws.Range("AC5").Value = 8
$ws.Range("AE5").Value = 21
$ws.Range("AK5").Value = 67
$ws.Range("AM5").Value = 51
$ws.Range("AT5").Value = 2.63
$ws.Range("AW5").Value = 7.5
$ws.Range("AX5").Value = 34

# Row 6
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 2.4

# Row 9
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 5.8
$ws.Range("J9").Value = 2.15
$ws.Range("L9").Value = 5.9
$ws.Range("N9").Value = 6.45
$ws.Range("O9").Value = 1.37
$ws.Range("P9").Value = 2.65
$ws.Range("T9").Value = 2.45
$ws.Range("W9").Value = 5.7
$ws.Range("X9").Value = 7
$ws.Range("Y9").Value = 8
$ws.Range("Z9").Value = 12.5
$ws.Range("AA9").Value = 14
$ws.Range("AC9").Value = 7.6
$ws.Range("AD9").Value = 6.6
$ws.Range("AE9").Value = 18.5
$ws.Range("AF9").Value = 110
$ws.Range("AG9").Value = 900
$ws.Range("AH9").Value = 12
$ws.Range("AI9").Value = 32
$ws.Range("AJ9").Value = 19
$ws.Range("AL9").Value = 80
$ws.Range("AM9").Value = 80
$ws.Range("AO9").Value = 7.7
$ws.Range("AP9").Value = 17.5
$ws.Range("AQ9").Value = 26
$ws.Range("AR9").Value = 55
$ws.Range("AU9").Value = 7.7
$ws.Range("AV9").Value = 80
$ws.Range("AX9").Value = 37
$ws.Range("AY9").Value = 40
$ws.Range("BA9").Value = 300
$ws.Range("BB9").Value = 450
